# Summary by Country and Scope
#
# This edit does three things to the "media_coverage" sheet:
#
#  1. Normalizes every multi-value Topic (column D) by removing the
#     space that follows each semicolon separator, e.g.
#     "Archaeology; History" -> "Archaeology;History".
#
#  2. Fills in a missing "Scope" (column B) value of
#     "International/National" for a handful of rows that had it blank.
#
#  3. Fixes a shifted row (ELLE Decor, row 177) where the Country value
#     ("Spain") was missing and "Architecture" had been entered in the
#     Country column instead of the Topic column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Collapse "; " -> ";" for every Topic cell in column D -------------
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $v = $cell.Value2
    if ($v -ne $null -and $v -is [string] -and $v.Contains("; ")) {
        $cell.Value = $v.Replace("; ", ";")
    }
}

# --- 2. Backfill missing Scope values --------------------------------------
$missingScopeRows = @(423, 429, 436, 445, 446, 456, 470, 476, 479, 489, 491, 494, 511, 514, 519, 529, 530, 533, 534, 553)
foreach ($r in $missingScopeRows) {
    $ws.Cells.Item($r, 2).Value = "International/National"
}

# --- 3. Fix row 177 (ELLE Decor): Country was missing, Topic had leaked ----
#        into the Country column.
$ws.Range("C177").Value = "Spain"
$ws.Range("D177").Value = "Architecture"
